$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume (and Coin/Link) columns in this sheet are stored as
# plain text, not numbers (e.g. "1.000", "30.218.86", "  -2.25%  "). Set
# each target cell to Text format before writing so Excel does not
# auto-convert the new value into a number/date and mangle it.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.229.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.878.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4839"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2870"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06580"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.901.08"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.77"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07322"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.127"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.07"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6543"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.204.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.32"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007743"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.387"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.125.80"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9993"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "192.90"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.122"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.253"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.78"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.01"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.910"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.433"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.258"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09079"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.007"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7145"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.098"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.700"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01777"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.636"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9225"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.72"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4267"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.785"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9991"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.383"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.80%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.73"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.825"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05757"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.70"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3810"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.31%  "
